# BIS-1002: remove the "Internal Assignment" column (column O) from the
# sample-type export: clear its header label and the per-row boolean
# values. The cells remain in place (retaining their existing style) but
# no longer carry any content, so the now-unused "Internal Assignment"
# shared string drops out of the workbook when it is saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4:O7").ClearContents()
